# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 48: The Sting of Conscience / Sleeping Potion
$ws.Range("H48").Value = 1356.75
$ws.Range("I48").Value = 809
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 2427
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = -2135
$ws.Range("N48").Value = -9584

# Row 56: Sleepless in Silvertear / Potent Sleeping Potion
$ws.Range("H56").Value = 1356.75
$ws.Range("I56").Value = 809
$ws.Range("J56").Value = 3000
$ws.Range("K56").Value = 2427
$ws.Range("L56").Value = 9000
$ws.Range("M56").Value = -1893
$ws.Range("N56").Value = -10068

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 784.6
$ws.Range("I80").Value = 817.41174
$ws.Range("K80").Value = 2452.23522
$ws.Range("M80").Value = -1454.23522

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 784.6
$ws.Range("I83").Value = 817.41174
$ws.Range("K83").Value = 7356.70566
$ws.Range("M83").Value = -2364.70566

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 1079.7
$ws.Range("I88").Value = 866.3333
$ws.Range("J88").Value = 1171.1428
$ws.Range("K88").Value = 866.3333
$ws.Range("L88").Value = 1171.1428
$ws.Range("M88").Value = -460.3333
$ws.Range("N88").Value = -1983.1428

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 1079.7
$ws.Range("I91").Value = 866.3333
$ws.Range("J91").Value = 1171.1428
$ws.Range("K91").Value = 866.3333
$ws.Range("L91").Value = 1171.1428
$ws.Range("M91").Value = 537.6667
$ws.Range("N91").Value = -3979.1428

# Row 121: Mindful Medicine / Tincture of Mind
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2207.2104
$ws.Range("I132").Value = 2012.2667
$ws.Range("K132").Value = 6036.800099999999
$ws.Range("M132").Value = -3506.800099999999

$ws = $wb.Worksheets.Item("ARM")
# Row 7: Distill It Yourself / Bronze Alembic
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2473.87
$ws.Range("I32").Value = 2494.818
$ws.Range("K32").Value = 2494.818
$ws.Range("M32").Value = -2207.818

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2074.1667
$ws.Range("I45").Value = 2236.25
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 2236.25
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -1859.25
$ws.Range("N45").Value = -2504

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3332.2
$ws.Range("I61").Value = 3332.2
$ws.Range("K61").Value = 3332.2
$ws.Range("M61").Value = -3120.2

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4790.5
$ws.Range("I122").Value = 3412.2
$ws.Range("K122").Value = 10236.6
$ws.Range("M122").Value = -7786.599999999999

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3175.475
$ws.Range("I132").Value = 3081.6216
$ws.Range("K132").Value = 9244.864799999999
$ws.Range("M132").Value = -6714.864799999999

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3332.2
$ws.Range("I136").Value = 3332.2
$ws.Range("K136").Value = 9996.599999999999
$ws.Range("M136").Value = -7446.599999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 3102.353
$ws.Range("I20").Value = 2853.4348
$ws.Range("K20").Value = 2853.4348
$ws.Range("M20").Value = -2606.4348

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 151228.58
$ws.Range("I134").Value = 9766.666999999999
$ws.Range("K134").Value = 29300.001
$ws.Range("M134").Value = -26765.001

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3326.375

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 2477.8
$ws.Range("I62").Value = 1989
$ws.Range("J62").Value = 2600
$ws.Range("K62").Value = 1989
$ws.Range("L62").Value = 2600
$ws.Range("M62").Value = -1365
$ws.Range("N62").Value = -3848

# Row 64: Almost as Fun as Slingshotting Birds / Cedar Longbow
$ws.Range("H64").Value = 110000
$ws.Range("J64").Value = 110000
$ws.Range("L64").Value = 110000
$ws.Range("N64").Value = -110496

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2477.8
$ws.Range("I65").Value = 1989
$ws.Range("J65").Value = 2600
$ws.Range("K65").Value = 9945
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = -6825
$ws.Range("N65").Value = -19240

# Row 67: Living Bow to Mouth (L) / Cedar Longbow
$ws.Range("H67").Value = 110000
$ws.Range("J67").Value = 110000
$ws.Range("L67").Value = 110000
$ws.Range("N67").Value = -111716

# Row 96: Composition / Larch Composite Bow
$ws.Range("H96").Value = 41966.668
$ws.Range("J96").Value = 41966.668
$ws.Range("L96").Value = 41966.668
$ws.Range("N96").Value = -47458.668

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2153.1333
$ws.Range("J132").Value = 3047.0908
$ws.Range("L132").Value = 9141.2724
$ws.Range("N132").Value = -14201.2724

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 458092.38
$ws.Range("I134").Value = 3648.15
$ws.Range("K134").Value = 10944.45
$ws.Range("M134").Value = -8409.450000000001

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3326.375

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 281.2
$ws.Range("I2").Value = 75.5
$ws.Range("J2").Value = 356
$ws.Range("K2").Value = 453
$ws.Range("L2").Value = 2136
$ws.Range("M2").Value = -340
$ws.Range("N2").Value = -2362

# Row 59: Comfort Me with Mushrooms / Buttons in a Blanket
$ws.Range("H59").Value = 150
$ws.Range("I59").Value = 150
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 450
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 90
$ws.Range("N59").ClearContents()

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 96474.55
$ws.Range("I107").Value = 1170.5
$ws.Range("J107").Value = 210839.4
$ws.Range("K107").Value = 3511.5
$ws.Range("L107").Value = 632518.2
$ws.Range("M107").Value = -1591.5
$ws.Range("N107").Value = -636358.2

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 33601.87
$ws.Range("I122").Value = 804
$ws.Range("J122").Value = 49219.906
$ws.Range("K122").Value = 7236
$ws.Range("L122").Value = 442979.154
$ws.Range("M122").Value = -4786
$ws.Range("N122").Value = -447879.154

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 7747.6
$ws.Range("I139").Value = 5107
$ws.Range("K139").Value = 15321
$ws.Range("M139").Value = -10181

$ws = $wb.Worksheets.Item("GSM")
# Row 38: He Has His Quartz / Silver Circlet (Goshenite)
$ws.Range("H38").Value = 59000
$ws.Range("J38").Value = 59000
$ws.Range("L38").Value = 59000
$ws.Range("N38").Value = -59926

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 6798.643
$ws.Range("I70").Value = 5477.263
$ws.Range("J70").Value = 9588.223
$ws.Range("K70").Value = 5477.263
$ws.Range("L70").Value = 9588.223
$ws.Range("M70").Value = -5207.263
$ws.Range("N70").Value = -10128.223

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 6798.643
$ws.Range("I73").Value = 5477.263
$ws.Range("J73").Value = 9588.223
$ws.Range("K73").Value = 5477.263
$ws.Range("L73").Value = 9588.223
$ws.Range("M73").Value = -4541.263
$ws.Range("N73").Value = -11460.223

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 3951.7896
$ws.Range("I102").Value = 1980.9166
$ws.Range("K102").Value = 1980.9166
$ws.Range("M102").Value = -358.9166

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3526.3635
$ws.Range("I126").Value = 1400
$ws.Range("K126").Value = 4200
$ws.Range("M126").Value = -1730

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1653
$ws.Range("I22").Value = 1663.8
$ws.Range("K22").Value = 1663.8
$ws.Range("M22").Value = -1368.8

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1653
$ws.Range("I27").Value = 1663.8
$ws.Range("K27").Value = 1663.8
$ws.Range("M27").Value = -1556.8

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 1722.2
$ws.Range("I55").Value = 62.166668
$ws.Range("K55").Value = 62.166668
$ws.Range("M55").Value = 110.833332

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2683.25
$ws.Range("I68").Value = 2592
$ws.Range("K68").Value = 2592
$ws.Range("M68").Value = -1843

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2683.25
$ws.Range("I71").Value = 2592
$ws.Range("K71").Value = 12960
$ws.Range("M71").Value = -9216

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 2199.4
$ws.Range("I82").Value = 2249.5
$ws.Range("J82").Value = 1999
$ws.Range("K82").Value = 2249.5
$ws.Range("L82").Value = 1999
$ws.Range("M82").Value = -1888.5
$ws.Range("N82").Value = -2721

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 2199.4
$ws.Range("I85").Value = 2249.5
$ws.Range("J85").Value = 1999
$ws.Range("K85").Value = 2249.5
$ws.Range("L85").Value = 1999
$ws.Range("M85").Value = -1001.5
$ws.Range("N85").Value = -4495

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 408.5
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 998
$ws.Range("I126").Value = 998
$ws.Range("K126").Value = 2994
$ws.Range("M126").Value = -524

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 29533.744
$ws.Range("I132").Value = 3060.3667
$ws.Range("J132").Value = 117778.336
$ws.Range("K132").Value = 9181.1001
$ws.Range("L132").Value = 353335.008
$ws.Range("M132").Value = -6651.1001
$ws.Range("N132").Value = -358395.008
